# Update the default "welcome" chatbot message and reflect the resulting
# UI state (active selection + row auto-height growth caused by the
# longer, wrapped text in B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New welcome message text (B2), replacing the old greeting.
$ws.Range("B2").Value = "Hello, my name is SAM. `nI am a virtual assistant, here to help you with information about Parkinson's Disease and Nuplazid."

# The cell that was last edited/clicked becomes the active selection.
$ws.Range("B2").Select()

# Row 2 grows taller to fit the longer wrapped message text.
$ws.Rows.Item(2).RowHeight = 45
